$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new blank rows before row 655 (shifts former rows 655-665 down to 661-671)
$ws.Range("A655:A660").EntireRow.Insert()

# Common (constant) values shared by all rows in this block
$mercadoId = 9
$mercado   = "Vega Central Mapocho de Santiago"
$region    = "Metropolitana"
$codreg    = 13
$tipo      = "Fruta"
$productoId = 100103
$producto   = "Frutos de hueso (carozo)"
$categoriaId = 100103006
$categoria   = "Nectarín"

# New rows data: Fecha, Variedad, Calidad, Volumen, PrecioMin, PrecioMax, PrecioProm, Unidad, Origen, PrecioKg, KgUnidad
$rows = @(
    @{ R=655; D=44628; K="Artic Snow";  L="Especial"; M=280; N=11200; O=11200; P=11200; Q="$/caja 16 kilos granel"; Rg="Región de O'Higgins";   S=700; T=16 },
    @{ R=656; D=44628; K="Artic Snow";  L="Primera";  M=330; N=9600;  O=9600;  P=9600;  Q="$/caja 16 kilos granel"; Rg="Región de O'Higgins";   S=600; T=16 },
    @{ R=657; D=44628; K="Artic Snow";  L="Segunda";  M=380; N=8000;  O=8000;  P=8000;  Q="$/caja 16 kilos granel"; Rg="Región de O'Higgins";   S=500; T=16 },
    @{ R=658; D=44628; K="August Red";  L="Especial"; M=280; N=11200; O=11200; P=11200; Q="$/caja 16 kilos granel"; Rg="Región Metropolitana";  S=700; T=16 },
    @{ R=659; D=44628; K="August Red";  L="Primera";  M=310; N=9600;  O=9600;  P=9600;  Q="$/caja 16 kilos granel"; Rg="Región Metropolitana";  S=600; T=16 },
    @{ R=660; D=44628; K="August Red";  L="Segunda";  M=350; N=8000;  O=8000;  P=8000;  Q="$/caja 16 kilos granel"; Rg="Región Metropolitana";  S=500; T=16 }
)

foreach ($row in $rows) {
    $r = $row.R
    $ws.Cells.Item($r, 1).Value  = $mercadoId
    $ws.Cells.Item($r, 2).Value  = $mercado
    $ws.Cells.Item($r, 3).Value  = $region
    $ws.Cells.Item($r, 4).Value  = $row.D
    $ws.Cells.Item($r, 5).Value  = $codreg
    $ws.Cells.Item($r, 6).Value  = $tipo
    $ws.Cells.Item($r, 7).Value  = $productoId
    $ws.Cells.Item($r, 8).Value  = $producto
    $ws.Cells.Item($r, 9).Value  = $categoriaId
    $ws.Cells.Item($r, 10).Value = $categoria
    $ws.Cells.Item($r, 11).Value = $row.K
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
    $ws.Cells.Item($r, 14).Value = $row.N
    $ws.Cells.Item($r, 15).Value = $row.O
    $ws.Cells.Item($r, 16).Value = $row.P
    $ws.Cells.Item($r, 17).Value = $row.Q
    $ws.Cells.Item($r, 18).Value = $row.Rg
    $ws.Cells.Item($r, 19).Value = $row.S
    $ws.Cells.Item($r, 20).Value = $row.T
}

Write-Output "done"
